$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "Wins", "Losses", "Ties" in AD1:AF1 ---
# Copy the format of the existing last header cell (AC1) so the new
# headers match the existing bold/centered/bordered header style,
# then overwrite with the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)
$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-46: season record (Wins=76, Losses=86, Ties=0) ---
$ws.Range("AD2:AD46").Value = 76
$ws.Range("AE2:AE46").Value = 86
$ws.Range("AF2:AF46").Value = 0

Write-Output "done"
